# Update TPM-derived LR-pair statistics (Il17c-Il17re) with new values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> ECs
$ws.Cells.Item(2,5).Value = 3
$ws.Cells.Item(2,6).Value = 1
$ws.Cells.Item(2,7).Value = 0.643452
$ws.Cells.Item(2,8).Value = 1.930356
$ws.Cells.Item(2,9).Value = 0.2431136893481813
$ws.Cells.Item(2,10).Value = 0.2431136893481813
$ws.Cells.Item(2,11).Value = 2
$ws.Cells.Item(2,12).Value = 0.6666666666666666
$ws.Cells.Item(2,13).Value = 0.03159966666666666
$ws.Cells.Item(2,14).Value = 0.094799
$ws.Cells.Item(2,15).Value = 0.07659813431479094
$ws.Cells.Item(2,16).Value = 0.07659813431479096
$ws.Cells.Item(2,17).Value = 0.020332868716
$ws.Cells.Item(2,18).Value = 0.182995818444
$ws.Cells.Item(2,19).Value = 0.01862205503045635
$ws.Cells.Item(2,20).Value = 0.01862205503045635

# Row 3: ECs -> FAPs
$ws.Cells.Item(3,5).Value = 3
$ws.Cells.Item(3,6).Value = 1
$ws.Cells.Item(3,7).Value = 0.643452
$ws.Cells.Item(3,8).Value = 1.930356
$ws.Cells.Item(3,9).Value = 0.2431136893481813
$ws.Cells.Item(3,10).Value = 0.2431136893481813
$ws.Cells.Item(3,15).Value = 0.428304440395438
$ws.Cells.Item(3,16).Value = 0.428304440395438
$ws.Cells.Item(3,17).Value = 0.113692820784
$ws.Cells.Item(3,18).Value = 1.023235387056
$ws.Cells.Item(3,19).Value = 0.1041266726687431
$ws.Cells.Item(3,20).Value = 0.1041266726687431

# Row 4: ECs -> MuSCs
$ws.Cells.Item(4,5).Value = 3
$ws.Cells.Item(4,6).Value = 1
$ws.Cells.Item(4,7).Value = 0.643452
$ws.Cells.Item(4,8).Value = 1.930356
$ws.Cells.Item(4,9).Value = 0.2431136893481813
$ws.Cells.Item(4,10).Value = 0.2431136893481813
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 0.171462
$ws.Cells.Item(4,14).Value = 0.514386
$ws.Cells.Item(4,15).Value = 0.4156268306379609
$ws.Cells.Item(4,16).Value = 0.415626830637961
$ws.Cells.Item(4,17).Value = 0.110327566824
$ws.Cells.Item(4,18).Value = 0.992948101416
$ws.Cells.Item(4,19).Value = 0.1010445721884864
$ws.Cells.Item(4,20).Value = 0.1010445721884864

# Row 5: ECs -> Resolving-Mac
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 0.643452
$ws.Cells.Item(5,8).Value = 1.930356
$ws.Cells.Item(5,9).Value = 0.2431136893481813
$ws.Cells.Item(5,10).Value = 0.2431136893481813
$ws.Cells.Item(5,11).Value = 3
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 0.03278466666666666
$ws.Cells.Item(5,14).Value = 0.098354
$ws.Cells.Item(5,15).Value = 0.07947059465181013
$ws.Cells.Item(5,16).Value = 0.07947059465181014
$ws.Cells.Item(5,17).Value = 0.021095359336
$ws.Cells.Item(5,18).Value = 0.189858234024
$ws.Cells.Item(5,19).Value = 0.0193203894604954
$ws.Cells.Item(5,20).Value = 0.01932038946049541

# Row 6: FAPs -> ECs
$ws.Cells.Item(6,9).Value = 0.2185191514957488
$ws.Cells.Item(6,10).Value = 0.2185191514957488
$ws.Cells.Item(6,11).Value = 2
$ws.Cells.Item(6,12).Value = 0.6666666666666666
$ws.Cells.Item(6,13).Value = 0.03159966666666666
$ws.Cells.Item(6,14).Value = 0.094799
$ws.Cells.Item(6,15).Value = 0.07659813431479094
$ws.Cells.Item(6,16).Value = 0.07659813431479096
$ws.Cells.Item(6,17).Value = 0.01827589894755556
$ws.Cells.Item(6,18).Value = 0.164483090528
$ws.Cells.Item(6,19).Value = 0.01673815931662551
$ws.Cells.Item(6,20).Value = 0.01673815931662552

# Row 7: FAPs -> FAPs
$ws.Cells.Item(7,9).Value = 0.2185191514957488
$ws.Cells.Item(7,10).Value = 0.2185191514957488
$ws.Cells.Item(7,15).Value = 0.428304440395438
$ws.Cells.Item(7,16).Value = 0.428304440395438
$ws.Cells.Item(7,19).Value = 0.09359272289707261
$ws.Cells.Item(7,20).Value = 0.09359272289707263

# Row 8: FAPs -> MuSCs
$ws.Cells.Item(8,9).Value = 0.2185191514957488
$ws.Cells.Item(8,10).Value = 0.2185191514957488
$ws.Cells.Item(8,11).Value = 3
$ws.Cells.Item(8,12).Value = 1
$ws.Cells.Item(8,13).Value = 0.171462
$ws.Cells.Item(8,14).Value = 0.514386
$ws.Cells.Item(8,15).Value = 0.4156268306379609
$ws.Cells.Item(8,16).Value = 0.415626830637961
$ws.Cells.Item(8,17).Value = 0.099166305088
$ws.Cells.Item(8,18).Value = 0.8924967457920001
$ws.Cells.Item(8,19).Value = 0.09082242236987449
$ws.Cells.Item(8,20).Value = 0.0908224223698745

# Row 9: FAPs -> Resolving-Mac
$ws.Cells.Item(9,9).Value = 0.2185191514957488
$ws.Cells.Item(9,10).Value = 0.2185191514957488
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 0.03278466666666666
$ws.Cells.Item(9,14).Value = 0.098354
$ws.Cells.Item(9,15).Value = 0.07947059465181013
$ws.Cells.Item(9,16).Value = 0.07947059465181014
$ws.Cells.Item(9,17).Value = 0.01896125238755555
$ws.Cells.Item(9,18).Value = 0.170651271488
$ws.Cells.Item(9,19).Value = 0.01736584691217614
$ws.Cells.Item(9,20).Value = 0.01736584691217614

# Row 10: MuSCs -> ECs
$ws.Cells.Item(10,7).Value = 0.2588786666666666
$ws.Cells.Item(10,8).Value = 0.776636
$ws.Cells.Item(10,9).Value = 0.09781141055878506
$ws.Cells.Item(10,10).Value = 0.09781141055878506
$ws.Cells.Item(10,11).Value = 2
$ws.Cells.Item(10,12).Value = 0.6666666666666666
$ws.Cells.Item(10,13).Value = 0.03159966666666666
$ws.Cells.Item(10,14).Value = 0.094799
$ws.Cells.Item(10,15).Value = 0.07659813431479094
$ws.Cells.Item(10,16).Value = 0.07659813431479096
$ws.Cells.Item(10,17).Value = 0.008180479573777776
$ws.Cells.Item(10,18).Value = 0.07362431616399999
$ws.Cells.Item(10,19).Value = 0.007492171563500979
$ws.Cells.Item(10,20).Value = 0.00749217156350098

# Row 11: MuSCs -> FAPs
$ws.Cells.Item(11,7).Value = 0.2588786666666666
$ws.Cells.Item(11,8).Value = 0.776636
$ws.Cells.Item(11,9).Value = 0.09781141055878506
$ws.Cells.Item(11,10).Value = 0.09781141055878506
$ws.Cells.Item(11,15).Value = 0.428304440395438
$ws.Cells.Item(11,16).Value = 0.428304440395438
$ws.Cells.Item(11,17).Value = 0.04574178937066666
$ws.Cells.Item(11,18).Value = 0.411676104336
$ws.Cells.Item(11,19).Value = 0.04189306146366887
$ws.Cells.Item(11,20).Value = 0.04189306146366888

# Row 12: MuSCs -> MuSCs
$ws.Cells.Item(12,7).Value = 0.2588786666666666
$ws.Cells.Item(12,8).Value = 0.776636
$ws.Cells.Item(12,9).Value = 0.09781141055878506
$ws.Cells.Item(12,10).Value = 0.09781141055878506
$ws.Cells.Item(12,11).Value = 3
$ws.Cells.Item(12,12).Value = 1
$ws.Cells.Item(12,13).Value = 0.171462
$ws.Cells.Item(12,14).Value = 0.514386
$ws.Cells.Item(12,15).Value = 0.4156268306379609
$ws.Cells.Item(12,16).Value = 0.415626830637961
$ws.Cells.Item(12,17).Value = 0.04438785394399999
$ws.Cells.Item(12,18).Value = 0.399490685496
$ws.Cells.Item(12,19).Value = 0.04065304657077622
$ws.Cells.Item(12,20).Value = 0.04065304657077622

# Row 13: MuSCs -> Resolving-Mac
$ws.Cells.Item(13,7).Value = 0.2588786666666666
$ws.Cells.Item(13,8).Value = 0.776636
$ws.Cells.Item(13,9).Value = 0.09781141055878506
$ws.Cells.Item(13,10).Value = 0.09781141055878506
$ws.Cells.Item(13,11).Value = 3
$ws.Cells.Item(13,12).Value = 1
$ws.Cells.Item(13,13).Value = 0.03278466666666666
$ws.Cells.Item(13,14).Value = 0.098354
$ws.Cells.Item(13,15).Value = 0.07947059465181013
$ws.Cells.Item(13,16).Value = 0.07947059465181014
$ws.Cells.Item(13,17).Value = 0.008487250793777776
$ws.Cells.Item(13,18).Value = 0.07638525714399999
$ws.Cells.Item(13,19).Value = 0.007773130960838989
$ws.Cells.Item(13,20).Value = 0.00777313096083899

# Row 14: Resolving-Mac -> ECs
$ws.Cells.Item(14,7).Value = 1.166024333333333
$ws.Cells.Item(14,8).Value = 3.498073
$ws.Cells.Item(14,9).Value = 0.4405557485972849
$ws.Cells.Item(14,10).Value = 0.4405557485972849
$ws.Cells.Item(14,11).Value = 2
$ws.Cells.Item(14,12).Value = 0.6666666666666666
$ws.Cells.Item(14,13).Value = 0.03159966666666666
$ws.Cells.Item(14,14).Value = 0.094799
$ws.Cells.Item(14,15).Value = 0.07659813431479094
$ws.Cells.Item(14,16).Value = 0.07659813431479096
$ws.Cells.Item(14,17).Value = 0.03684598025855555
$ws.Cells.Item(14,18).Value = 0.331613822327
$ws.Cells.Item(14,19).Value = 0.0337457484042081
$ws.Cells.Item(14,20).Value = 0.0337457484042081

# Row 15: Resolving-Mac -> FAPs
$ws.Cells.Item(15,7).Value = 1.166024333333333
$ws.Cells.Item(15,8).Value = 3.498073
$ws.Cells.Item(15,9).Value = 0.4405557485972849
$ws.Cells.Item(15,10).Value = 0.4405557485972849
$ws.Cells.Item(15,15).Value = 0.428304440395438
$ws.Cells.Item(15,16).Value = 0.428304440395438
$ws.Cells.Item(15,17).Value = 0.2060271715053333
$ws.Cells.Item(15,18).Value = 1.854244543548
$ws.Cells.Item(15,19).Value = 0.1886919833659534
$ws.Cells.Item(15,20).Value = 0.1886919833659534

# Row 16: Resolving-Mac -> MuSCs
$ws.Cells.Item(16,7).Value = 1.166024333333333
$ws.Cells.Item(16,8).Value = 3.498073
$ws.Cells.Item(16,9).Value = 0.4405557485972849
$ws.Cells.Item(16,10).Value = 0.4405557485972849
$ws.Cells.Item(16,11).Value = 3
$ws.Cells.Item(16,12).Value = 1
$ws.Cells.Item(16,13).Value = 0.171462
$ws.Cells.Item(16,14).Value = 0.514386
$ws.Cells.Item(16,15).Value = 0.4156268306379609
$ws.Cells.Item(16,16).Value = 0.415626830637961
$ws.Cells.Item(16,17).Value = 0.199928864242
$ws.Cells.Item(16,18).Value = 1.799359778178
$ws.Cells.Item(16,19).Value = 0.1831067895088238
$ws.Cells.Item(16,20).Value = 0.1831067895088238

# Row 17: Resolving-Mac -> Resolving-Mac
$ws.Cells.Item(17,7).Value = 1.166024333333333
$ws.Cells.Item(17,8).Value = 3.498073
$ws.Cells.Item(17,9).Value = 0.4405557485972849
$ws.Cells.Item(17,10).Value = 0.4405557485972849
$ws.Cells.Item(17,11).Value = 3
$ws.Cells.Item(17,12).Value = 1
$ws.Cells.Item(17,13).Value = 0.03278466666666666
$ws.Cells.Item(17,14).Value = 0.098354
$ws.Cells.Item(17,15).Value = 0.07947059465181013
$ws.Cells.Item(17,16).Value = 0.07947059465181014
$ws.Cells.Item(17,17).Value = 0.03822771909355555
$ws.Cells.Item(17,18).Value = 0.344049471842
$ws.Cells.Item(17,19).Value = 0.03501122731829959
$ws.Cells.Item(17,20).Value = 0.0350112273182996

